# Apply crypto price/volume updates per commit "Updated cryptos list on Sun Oct 15 11:14:25 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.074.31"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").Value = "'1.569.52"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("E4").Value = "  +0.78%  "
$ws.Range("D5").Value = "'209.15"
$ws.Range("E5").Value = "  +1.39%  "
$ws.Range("D6").Value = "'0.492"
$ws.Range("E6").Value = "  +0.85%  "
$ws.Range("E7").Value = "  +0.65%  "
$ws.Range("D8").Value = "'22.14"
$ws.Range("E8").Value = "  +0.47%  "
$ws.Range("D9").Value = "'0.250"
$ws.Range("E9").Value = "  +1.13%  "
$ws.Range("E10").Value = "  +1.68%  "
$ws.Range("D11").Value = "'0.0862"
$ws.Range("E11").Value = "  +0.87%  "
$ws.Range("D12").Value = "'1.576.06"
$ws.Range("E12").Value = "  +2.08%  "
$ws.Range("D13").Value = "'3.78"
$ws.Range("E13").Value = "  +1.18%  "
$ws.Range("D14").Value = "'0.521"
$ws.Range("E14").Value = "  +0.38%  "
$ws.Range("D15").Value = "'27.070.73"
$ws.Range("E15").Value = "  +0.58%  "
$ws.Range("D16").Value = "'62.04"
$ws.Range("E16").Value = "  +0.60%  "
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").Value = "'216.10"
$ws.Range("E18").Value = "  -0.57%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'7.42"
$ws.Range("E19").Value = "  +2.28%  "
$ws.Range("E20").Value = "  +0.70%  "
$ws.Range("D21").Value = "'4.16"
$ws.Range("E21").Value = "  +2.53%  "
$ws.Range("D22").Value = "'9.20"
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("D23").Value = "'1.95"
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("D24").Value = "'154.35"
$ws.Range("E24").Value = "  +0.36%  "
$ws.Range("D25").Value = "'6.64"
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").Value = "'15.07"
$ws.Range("E26").Value = "  +0.81%  "
$ws.Range("D27").Value = "'0.106"
$ws.Range("E27").Value = "  +1.43%  "
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "'1.13"
$ws.Range("E29").Value = "  +4.24%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "'0.0474"
$ws.Range("E30").Value = "  +1.21%  "
$ws.Range("D31").Value = "'3.24"
$ws.Range("E31").Value = "  +0.50%  "
$ws.Range("E32").Value = "  +3.04%  "
$ws.Range("D33").Value = "'1.427.20"
$ws.Range("E33").Value = "  +1.36%  "
$ws.Range("E34").Value = "  +12.88%  "
$ws.Range("E35").Value = "  +1.08%  "
$ws.Range("D36").Value = "'2.37"
$ws.Range("E36").Value = "  +3.78%  "
$ws.Range("D37").Value = "'0.0168"
$ws.Range("E37").Value = "  +1.28%  "
$ws.Range("D38").Value = "'0.534"
$ws.Range("E38").Value = "  +1.04%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.45"
$ws.Range("E39").Value = "  +6.36%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'5.84"
$ws.Range("E40").Value = "  +2.75%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "'0.814"
$ws.Range("E41").Value = "  +1.03%  "
$ws.Range("E42").Value = "  +0.87%  "
$ws.Range("E43").Value = "  +1.17%  "
$ws.Range("D44").Value = "'64.78"
$ws.Range("E44").Value = "  +0.44%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").Value = "'1.706.67"
$ws.Range("E46").Value = "  +1.24%  "
$ws.Range("D47").Value = "'86.63"
$ws.Range("E47").Value = "  -0.91%  "
$ws.Range("D48").Value = "'0.0₆0102"
$ws.Range("E48").Value = "  +1.81%  "
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("D50").Value = "'0.0965"
$ws.Range("E50").Value = "  +0.48%  "
$ws.Range("E51").Value = "  +0.62%  "
